$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '37.080.25'
$ws.Cells.Item(2, 5).Value = '  -1.65%  '
$ws.Cells.Item(3, 4).Value = '2.020.90'
$ws.Cells.Item(3, 5).Value = '  -2.70%  '
$ws.Cells.Item(4, 5).Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '226.34'
$ws.Cells.Item(5, 5).Value = '  -2.70%  '
$ws.Cells.Item(6, 5).Value = '  -2.64%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '54.89'
$ws.Cells.Item(8, 5).Value = '  -5.37%  '
$ws.Cells.Item(9, 5).Value = '  -3.53%  '
$ws.Cells.Item(10, 5).Value = '  +0.19%  '
$ws.Cells.Item(11, 5).Value = '  -5.56%  '
$ws.Cells.Item(12, 4).Value = '2.317.69'
$ws.Cells.Item(12, 5).Value = '  -2.79%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '14.13'
$ws.Cells.Item(13, 5).Value = '  -5.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '20.24'
$ws.Cells.Item(14, 5).Value = '  -4.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.740'
$ws.Cells.Item(15, 5).Value = '  -3.37%  '
$ws.Cells.Item(16, 5).Value = '  -3.74%  '
$ws.Cells.Item(17, 4).Value = '2.016.69'
$ws.Cells.Item(17, 5).Value = '  -3.19%  '
$ws.Cells.Item(18, 4).Value = '37.045.62'
$ws.Cells.Item(18, 5).Value = '  -1.63%  '
$ws.Cells.Item(19, 5).Value = '  +0.96%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '68.96'
$ws.Cells.Item(20, 5).Value = '  -1.89%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0819'
$ws.Cells.Item(21, 5).Value = '  -1.60%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '223.56'
$ws.Cells.Item(22, 5).Value = '  -1.82%  '
$ws.Cells.Item(23, 5).Value = '  +0.06%  '
$ws.Cells.Item(24, 5).Value = '  +2.07%  '
$ws.Cells.Item(25, 5).Value = '  -7.67%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '166.18'
$ws.Cells.Item(26, 5).Value = '  -2.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '9.18'
$ws.Cells.Item(27, 5).Value = '  -7.88%  '
$ws.Cells.Item(28, 5).Value = '  -2.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '18.75'
$ws.Cells.Item(29, 5).Value = '  -3.11%  '
$ws.Cells.Item(30, 5).Value = '  -5.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.117'
$ws.Cells.Item(31, 5).Value = '  -3.83%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.50'
$ws.Cells.Item(32, 5).Value = '  -2.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.0613'
$ws.Cells.Item(33, 5).Value = '  -2.96%  '
$ws.Cells.Item(34, 5).Value = '  -5.32%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '2.35'
$ws.Cells.Item(35, 5).Value = '  -7.25%  '
$ws.Cells.Item(36, 5).Value = '  +1.51%  '
$ws.Cells.Item(37, 5).Value = '  -0.19%  '
$ws.Cells.Item(38, 5).Value = '  -5.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '5.28'
$ws.Cells.Item(39, 5).Value = '  -1.10%  '
$ws.Cells.Item(40, 2).Value = 'Maker'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(40, 4).Value = '1.474.96'
$ws.Cells.Item(40, 5).Value = '  -1.00%  '
$ws.Cells.Item(41, 2).Value = 'VeChain'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.0216'
$ws.Cells.Item(41, 5).Value = '  -5.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '95.45'
$ws.Cells.Item(42, 5).Value = '  -3.25%  '
$ws.Cells.Item(43, 5).Value = '  -4.31%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '16.28'
$ws.Cells.Item(44, 5).Value = '  -4.29%  '
$ws.Cells.Item(45, 5).Value = '  -5.26%  '
$ws.Cells.Item(46, 5).Value = '  -6.08%  '
$ws.Cells.Item(47, 5).Value = '  -3.25%  '
$ws.Cells.Item(48, 5).Value = '  -1.36%  '
$ws.Cells.Item(49, 5).Value = '  -1.54%  '
$ws.Cells.Item(50, 4).Value = '2.204.95'
$ws.Cells.Item(50, 5).Value = '  -2.82%  '
$ws.Cells.Item(51, 5).Value = '  -12.64%  '
